$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting rows 5-23 down to 6-24
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "regions_countries"
$ws.Range("B5").Value = "Change in regional GHGs from multiple perspectives and their underlying drivers"
$ws.Range("C5").Value = "Panel a: Regional GHG emission trends (in GtCO2eq yr-1) for the time period 1990-2019. GHG emissions from international aviation (AIR) and shipping (SEA) are not assigned to individual countries and shown separately. Panels b and c: Changes in GHG emissions for largest emitters (75% of global emissions) for the post-AR5 reporting period 2010-2019 in relative (% annual change) and absolute terms (GtCO2eq). Panels d and e: GHG emissions per capita and per GDP in 2019 for the largest emitters (75% of global emissions). GDP estimated using constant international purchasing power parity (USD 2017). Emissions are converted into CO2-equivalents based on global warming potentials with a 100 year time horizon (GWP-100) from the IPCC Sixth Assessment Report. The yellow dots represent the emissions data according to the national inventories reported by the Annex I countries of the Kyoto Protocol to the UNFCCC (Gütschow et al., 2021; Louise Jeffery et al., 2018).  LULUCF CO2 emissions are included in panel a, based on the average of three bookkeeping models (see data explanation in section 2.2 introduction), but are excluded in panels b to e due to a lack of regional resolution."

$ws.Range("C6").Select()
